{"js": "// Update the worksheet header date and each \"AxB=C\" answer cell to the\n// newly generated values. Every old value is unique in the document, so a\n// scoped body.search() + insertText(\"Replace\") per pair is unambiguous.\nconst replacements = [\n  [\"2026-02-16 Monday\", \"2026-02-17 Tuesday\"],\n  [\"80\u00d792=7360\", \"64\u00d779=5056\"],\n  [\"96\u00d737=3552\", \"88\u00d719=1672\"],\n  [\"46\u00d721=966\", \"67\u00d722=1474\"],\n  [\"25\u00d723=575\", \"26\u00d798=2548\"],\n  [\"73\u00d775=5475\", \"26\u00d798=2548\"],\n  [\"92\u00d760=5520\", \"24\u00d752=1248\"],\n  [\"28\u00d733=924\", \"24\u00d761=1464\"],\n  [\"55\u00d719=1045\", \"39\u00d756=2184\"],\n  [\"83\u00d723=1909\", \"88\u00d798=8624\"],\n  [\"74\u00d721=1554\", \"81\u00d735=2835\"],\n  [\"79\u00d724=1896\", \"40\u00d737=1480\"],\n  [\"95\u00d787=8265\", \"79\u00d799=7821\"],\n  [\"19\u00d728=532\", \"69\u00d772=4968\"],\n  [\"88\u00d782=7216\", \"50\u00d788=4400\"],\n  [\"46\u00d756=2576\", \"32\u00d744=1408\"],\n  [\"46\u00d797=4462\", \"38\u00d769=2622\"],\n  [\"50\u00d786=4300\", \"26\u00d785=2210\"],\n  [\"20\u00d734=680\", \"23\u00d738=874\"],\n  [\"93\u00d785=7905\", \"62\u00d783=5146\"],\n  [\"55\u00d732=1760\", \"77\u00d723=1771\"],\n  [\"45\u00d723=1035\", \"78\u00d713=1014\"],\n  [\"84\u00d771=5964\", \"42\u00d793=3906\"],\n  [\"59\u00d788=5192\", \"18\u00d746=828\"],\n  [\"86\u00d750=4300\", \"77\u00d741=3157\"],\n  [\"59\u00d728=1652\", \"72\u00d773=5256\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the worksheet header date and each \"AxB=C\" answer cell to the\n# newly generated values. Every old value is unique in the document, so a\n# single Find/Replace (wdReplaceOne semantics via a fresh Find scoped to\n# $d.Content each time) per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-16 Monday\", \"2026-02-17 Tuesday\"),\n    @(\"80\u00d792=7360\", \"64\u00d779=5056\"),\n    @(\"96\u00d737=3552\", \"88\u00d719=1672\"),\n    @(\"46\u00d721=966\", \"67\u00d722=1474\"),\n    @(\"25\u00d723=575\", \"26\u00d798=2548\"),\n    @(\"73\u00d775=5475\", \"26\u00d798=2548\"),\n    @(\"92\u00d760=5520\", \"24\u00d752=1248\"),\n    @(\"28\u00d733=924\", \"24\u00d761=1464\"),\n    @(\"55\u00d719=1045\", \"39\u00d756=2184\"),\n    @(\"83\u00d723=1909\", \"88\u00d798=8624\"),\n    @(\"74\u00d721=1554\", \"81\u00d735=2835\"),\n    @(\"79\u00d724=1896\", \"40\u00d737=1480\"),\n    @(\"95\u00d787=8265\", \"79\u00d799=7821\"),\n    @(\"19\u00d728=532\", \"69\u00d772=4968\"),\n    @(\"88\u00d782=7216\", \"50\u00d788=4400\"),\n    @(\"46\u00d756=2576\", \"32\u00d744=1408\"),\n    @(\"46\u00d797=4462\", \"38\u00d769=2622\"),\n    @(\"50\u00d786=4300\", \"26\u00d785=2210\"),\n    @(\"20\u00d734=680\", \"23\u00d738=874\"),\n    @(\"93\u00d785=7905\", \"62\u00d783=5146\"),\n    @(\"55\u00d732=1760\", \"77\u00d723=1771\"),\n    @(\"45\u00d723=1035\", \"78\u00d713=1014\"),\n    @(\"84\u00d771=5964\", \"42\u00d793=3906\"),\n    @(\"59\u00d788=5192\", \"18\u00d746=828\"),\n    @(\"86\u00d750=4300\", \"77\u00d741=3157\"),\n    @(\"59\u00d728=1652\", \"72\u00d773=5256\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n\n    # wdFindContinue=1, wdReplaceAll=2; MatchCase/MatchWholeWord true so each\n    # lookup is exact (values are unique, so this touches exactly one run).\n    $ok = $find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $ok) {\n        throw \"Replace failed for: $findText\"\n    }\n}\n"}
